$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the week-2 / week-3 values between row 2 and row 3
# (columns D, M, N, O, P, Q, S, T)

$ws.Range("D2").Value = 44973
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("Q2").Value = "$/bandeja 5 kilos"
$ws.Range("S2").Value = 2400
$ws.Range("T2").Value = 5

$ws.Range("D3").Value = 44238
$ws.Range("M3").Value = 35
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 2000
$ws.Range("T3").Value = 10
